# INRC2_Simulator/log/TestedParameter.xlsx
# Add a new logged run (column L) with the MIN_TABU_BASE parameter results,
# add a new note (F11/F12) about the ARBCS/ARRCS run, and drop the stale
# I12 rank entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: results of the newly logged test run ------------------
$ws.Range("L1").Value = "123"
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 0.2
$ws.Range("L4").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("L7").Value = 0.2
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("L13").Value = "0424"

# --- New note cell F11 (wraps across two lines) ---------------------------
$ws.Range("F11").Value = "ARBCS`nARRCS"
$ws.Range("F11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 27

# --- New rank sample F12, and drop the stale I12 rank ---------------------
$ws.Range("F12").Value = 713
$ws.Range("I12").Clear()

# --- Move the active selection in the frozen bottom-right pane ------------
$ws.Range("L15").Select()
